# Jasprit Bumrah.xlsx - full scrape update
# 1) Rename the sheet
# 2) Insert a new "matchNo" column at A (shifting the old columns right by one)
# 3) Replace the single data row with the 7 scraped rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Jasprit Bumrah"

# Header row
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($c = 1; $c -le $headers.Count; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows (matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr, opponentTeamName, venue, date, result)
$data = @(
    @("55th","Mumbai Indians","Jasprit Bumrah","","5","2","1","0","250.00","Sunrisers Hyderabad","Abu Dhabi","October 08","Mumbai won by 42 runs"),
    @("46th","Mumbai Indians","Jasprit Bumrah","","1","1","0","0","100.00","Delhi Capitals","Sharjah","October 02","Capitals won by 4 wickets (with 5 balls remaining)"),
    @("13th","Mumbai Indians","Jasprit Bumrah","","3","3","0","0","100.00","Delhi Capitals","Chennai","April 20","Capitals won by 6 wickets (with 5 balls remaining)"),
    @("5th","Mumbai Indians","Jasprit Bumrah","c Shakib Al Hasan b Russell","0","1","0","0","0.00","Kolkata Knight Riders","Chennai","April 13","Mumbai won by 10 runs"),
    @("30th","Mumbai Indians","Jasprit Bumrah","","1","2","0","0","50.00","Chennai Super Kings","Dubai (DSC)","September 19","Super Kings won by 20 runs"),
    @("39th","Mumbai Indians","Jasprit Bumrah","b Chahal","5","6","1","0","83.33","Royal Challengers Bangalore","Dubai (DSC)","September 26","RCB won by 54 runs"),
    @("1st","Mumbai Indians","Jasprit Bumrah","","1","2","0","0","50.00","Royal Challengers Bangalore","Chennai","April 09","RCB won by 2 wickets")
)

# Columns whose values are textual numbers and must stay text (not real numbers)
$numericLikeCols = @(5, 6, 7, 8, 9)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 1; $c -le $row.Count; $c++) {
        $val = $row[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($val -eq "") {
            # Keep an explicit empty-text cell (matches source "" <v/> cells)
            $cell.Formula = '=""'
        } elseif ($numericLikeCols -contains $c) {
            # Force storage as text so "5", "250.00", "83.33" etc. keep their
            # literal text representation instead of becoming real numbers
            $escaped = $val -replace '"', '""'
            $cell.Formula = '="' + $escaped + '"'
        } else {
            $cell.Value = $val
        }
    }
}
